$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header cells B1:P1 get ".jamais" suffix replaced with ".jamais.age_trait"
$headerCells = @("B1", "C1", "D1", "E1", "F1", "G1", "H1", "I1", "J1", "K1", "L1", "M1", "N1", "O1", "P1")

foreach ($addr in $headerCells) {
    $cell = $ws.Range($addr)
    $current = $cell.Value2
    if ($current -ne $null -and $current.ToString().EndsWith(".jamais")) {
        $cell.Value = $current.ToString() + ".age_trait"
    }
}
